# Fruta / hortaliza, semanal
# Update Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) for several rows
# of the "Arándano (blue)" sheet to reflect the corrected weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44463
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13500
$ws.Range("S2").Value = 6750

# Row 4
$ws.Range("D4").Value = 44497
$ws.Range("M4").Value = 400

# Row 5
$ws.Range("D5").Value = 44491
$ws.Range("M5").Value = 200

# Row 7
$ws.Range("D7").Value = 44468
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 13500
$ws.Range("S7").Value = 6750

# Row 8
$ws.Range("D8").Value = 44495
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11500
$ws.Range("S8").Value = 5750

# Row 9
$ws.Range("D9").Value = 44466
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 13500
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 13750
$ws.Range("S9").Value = 6875

# Row 10
$ws.Range("D10").Value = 44452
$ws.Range("M10").Value = 200

# Row 11
$ws.Range("D11").Value = 44494
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 11500
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 11750
$ws.Range("S11").Value = 5875

# Row 12
$ws.Range("D12").Value = 44454
$ws.Range("M12").Value = 300

# Row 13
$ws.Range("D13").Value = 44455
$ws.Range("M13").Value = 160

# Row 14
$ws.Range("D14").Value = 44490
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 11500
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 11750
$ws.Range("S14").Value = 5875

# Row 15
$ws.Range("D15").Value = 44462
$ws.Range("M15").Value = 140

# Row 16
$ws.Range("D16").Value = 44446
$ws.Range("M16").Value = 300

# Row 17
$ws.Range("D17").Value = 44445
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = 14000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 14500
$ws.Range("S17").Value = 7250

# Row 19
$ws.Range("D19").Value = 44498
$ws.Range("M19").Value = 240
$ws.Range("N19").Value = 11000
$ws.Range("O19").Value = 11500
$ws.Range("P19").Value = 11250
$ws.Range("S19").Value = 5625
